$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Only the cells whose value actually changes (runs/balls/fours/sixes
# columns C-F were reshuffled across the match rows for Shivam Dube).
$updates = @{
    "C2"  = "11"
    "D2"  = "12"
    "F2"  = "1"
    "C3"  = "23"
    "D3"  = "19"
    "E3"  = "0"
    "F3"  = "2"
    "C4"  = "27"
    "D4"  = "10"
    "E4"  = "1"
    "F4"  = "3"
    "C5"  = "17"
    "D5"  = "11"
    "E5"  = "2"
    "C6"  = "8"
    "D6"  = "13"
    "F6"  = "0"
    "C8"  = "7"
    "D8"  = "8"
    "E8"  = "0"
    "F8"  = "0"
    "C9"  = "12"
    "D9"  = "12"
    "F9"  = "1"
    "C10" = "2"
    "D10" = "6"
    "F10" = "0"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Prefix with an apostrophe so the numeric-looking text is stored as
    # text (matching the source data's string-typed cells) instead of
    # being auto-converted to a number, then restore the default style
    # so no stray quote-prefix formatting is left behind.
    $cell.Value = "'" + $updates[$addr]
    $cell.Style = "Normal"
}
